# Ejercicios 1 al 10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Duplicate the "Diametro" frequency table (A19:B28) into D40:E49, matching
#    the borders/number-format look of the original boxed mini-table.
# ---------------------------------------------------------------------------
$ws.Range("A19:B28").Copy()
$ws.Range("D40").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$diam = @(25.26, 25.27, 25.3, 25.31, 25.32, 25.33, 25.34, 25.35, 25.37, 25.38)
$freq = @(1, 1, 1, 1, 1, 1, 1, 2, 1, 2)
for ($i = 0; $i -lt 10; $i++) {
    $r = 40 + $i
    $ws.Cells.Item($r, 4).Value = $diam[$i]
    $ws.Cells.Item($r, 5).Value = $freq[$i]
}

# ---------------------------------------------------------------------------
# 2) New binned-range frequency tables in rows 53:57 (C:D descending by bin,
#    F:G ascending by bin -- F:G is what backs the updated chart).
# ---------------------------------------------------------------------------
$binsDesc = @("25,36-25,38", "25,33-25,35", "25,30-25,32", "25,29-25,31", "25,26-25,28")
$countsDesc = @(3, 4, 3, 2, 2)
$binsAsc = @("25,26-25,28", "25,29-25,31", "25,30-25,32", "25,33-25,35", "25,36-25,38")
$countsAsc = @(2, 2, 3, 4, 3)

for ($i = 0; $i -lt 5; $i++) {
    $r = 53 + $i
    $ws.Cells.Item($r, 3).Value = $binsDesc[$i]
    $ws.Cells.Item($r, 4).Value = $countsDesc[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $r = 53 + $i
    $ws.Cells.Item($r, 6).Value = $binsAsc[$i]
    $ws.Cells.Item($r, 7).Value = $countsAsc[$i]
}
$ws.Cells.Item(53, 8).Value = 2

# ---------------------------------------------------------------------------
# 3) Point the "Diametro" histogram (second chart) at the new binned table.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects(2)
$chart = $co.Chart
for ($i = 1; $i -le $chart.SeriesCollection().Count; $i++) {
    $ser = $chart.SeriesCollection($i)
    $ser.Values = $ws.Range("G53:G57")
    $ser.XValues = $ws.Range("F53:F57")
}

$catAxis = $chart.Axes(1)
$catAxis.TickLabels.NumberFormat = "General"

# ---------------------------------------------------------------------------
# 4) Reposition/resize the chart now that the sheet has grown.
# ---------------------------------------------------------------------------
$co.Left = $ws.Range("I27").Left
$co.Top = $ws.Range("I27").Top
$co.Width = 425
$co.Height = 310

# ---------------------------------------------------------------------------
# 5) View cosmetics to mirror the author's session.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 115
$ws.Range("S40").Select()
